# Updated symbol list on Sat Dec 24 13:27:41 UTC 2022 with GitHub Actions
# This script applies the latest price-refresh values scraped for the crypto tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "'244.75"
$ws.Range("D3").Value = "'21.90"
$ws.Range("D5").Value = "'0.05995"
$ws.Range("D6").Value = "'3.392"
$ws.Range("D7").Value = "'0.8142"
$ws.Range("D8").Value = "'0.9566"
$ws.Range("D9").Value = "'0.1425"
$ws.Range("D10").Value = "'0.07412"
$ws.Range("D11").Value = "'0.03334"
$ws.Range("D12").Value = "'0.03058"
$ws.Range("D13").Value = "'0.09408"
$ws.Range("D14").Value = "'4.003"
$ws.Range("D15").Value = "'0.001595"
$ws.Range("D16").Value = "'0.04810"
$ws.Range("D18").Value = "'0.006179"
$ws.Range("D19").Value = "'0.004999"
$ws.Range("D20").Value = "'0.0009881"
$ws.Range("D23").Value = "'6.416"
$ws.Range("D26").Value = "'0.1324"
$ws.Range("D27").Value = "'0.0002448"
$ws.Range("D40").Value = "'0.03988"
$ws.Range("D41").Value = "'0.006572"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D43").Value = "'0.002721"
$ws.Range("D44").Value = "'0.005764"
$ws.Range("D45").Value = "'0.00005278"
$ws.Range("D47").Value = "'1.000"
$ws.Range("D48").Value = "'0.01370"
$ws.Range("D49").Value = "'0.00002100"

# --- Column E (Volume(1h) label) updates ---
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
$ws.Range("E41").Value = "40KickTokenKICK"

